$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TotalCost (column C) values for the 96 scenarios where ID_PV = 1.
# New values reflect a refreshed operation-energy-cost calculation run.
$updates = @{
    2 = 120716.1522112524
    3 = 139051.7984401947
    6 = 122584.2638928925
    7 = 143467.8780872832
    10 = 121897.7867472824
    11 = 141028.3171938551
    14 = 124177.2568180901
    15 = 145824.7654788709
    18 = 142022.830156972
    19 = 170923.5954583264
    22 = 141864.4580373701
    23 = 170794.4238972236
    26 = 141874.069189175
    27 = 170804.409701393
    30 = 141715.9680147173
    31 = 170675.3451249267
    34 = 184772.8045993561
    35 = 200441.6653923719
    38 = 186148.520469406
    39 = 204072.0946239115
    42 = 186036.1053038605
    43 = 202931.8747747873
    46 = 187927.231350612
    47 = 207255.2241745333
    50 = 207211.4476988417
    51 = 236146.7467773402
    54 = 207053.0755731615
    55 = 236017.5756465657
    58 = 207062.6867280056
    59 = 236027.5612394723
    62 = 206904.5855535479
    63 = 235898.4966656377
    66 = 234906.859567781
    67 = 249470.2379004707
    70 = 236064.9238989371
    71 = 252786.2787378968
    74 = 236635.6795396588
    75 = 252773.824457818
    78 = 238447.614645082
    79 = 256959.3831505053
    82 = 260773.8126323012
    83 = 289718.9247643672
    86 = 260615.4405126992
    87 = 289589.7534118955
    90 = 260625.0516645042
    91 = 289599.7390148292
    94 = 260466.9504900465
    95 = 289470.6744309675
    98 = 125819.7983852314
    99 = 152295.34210942
    102 = 131660.3520990623
    103 = 158685.4376367659
    106 = 125819.7983852314
    107 = 152295.34210942
    110 = 131660.3520990623
    111 = 158685.4376367659
    114 = 142591.0276293001
    115 = 171514.5651693879
    118 = 141799.4095141054
    119 = 170675.3451249266
    122 = 142591.0276293001
    123 = 171514.5651693879
    126 = 141799.4095141054
    127 = 170675.3451249266
    130 = 190229.2404402616
    131 = 214131.6591120306
    134 = 195648.0359737862
    135 = 219586.4475331511
    138 = 190229.2404402616
    139 = 214131.6591120306
    142 = 195648.0359737862
    143 = 219586.4475331511
    146 = 207814.1791700112
    147 = 236737.716710099
    150 = 207022.5610548165
    151 = 235898.4966656376
    154 = 207814.1791700112
    155 = 236737.716710099
    158 = 207022.5610548165
    159 = 235898.4966656376
    162 = 240340.4944477063
    163 = 263076.5815397532
    166 = 245577.5501851852
    167 = 268186.0062417226
    170 = 240340.4944477063
    171 = 263076.5815397532
    174 = 245577.5501851852
    175 = 268186.0062417226
    178 = 261386.356935341
    179 = 290309.8944754287
    182 = 260594.7388201463
    183 = 289470.6744309674
    186 = 261386.356935341
    187 = 290309.8944754287
    190 = 260594.7388201463
    191 = 289470.6744309674
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

